# Add an "ORCID" column (D) with author ORCID ids for the first six
# authors, give it the same wrap-text formatting the sheet already uses,
# size the column, and normalize C15's font to match the other email
# hyperlink cells (C4/C6/C7/C8) instead of its separate duplicate font.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ORCID values for rows 1-6 (column D)
$orcids = @{
    1 = "0000-0001-5990-6585"
    2 = "0000-0002-7669-7364"
    3 = "0009-0000-8247-7432"
    4 = "0009-0007-7427-9437"
    5 = "0009-0006-5483-3667"
    6 = "0009-0000-0309-0353 "
}

foreach ($row in 1..6) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = $orcids[$row]
    $cell.WrapText = $true
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
}

# Size the new column like the rest of the author-info table
$ws.Columns.Item(4).ColumnWidth = 24.86

# C15 ("Zhang Shenli" e-mail) should use the same visual style as the
# other hyperlinked e-mail cells (C4, C6, C7, C8) rather than its own
# separate (duplicate) font definition.
$ws.Range("C4").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the author's final selection in the sheet
$ws.Range("F8").Select()
